# Filter CytokineExpression on SNP ID
# Add a new "snp" xref attribute pointing to a new SnpsToPlot EMX file,
# and add the Cytokine / Stimulus lookup sheets.

$wb = $excel.ActiveWorkbook
$attrs = $wb.Worksheets.Item(1)

# snp attribute (row 5) becomes an xref into the new SnpsToPlot entity
$attrs.Range("C5").Value = "xref"
$attrs.Range("F5").Value = "SnpsToPlot"

# Try to nudge the window position to mirror the author's saved view state.
try {
    $wb.Windows.Item(1).Left = 7880
} catch {
}

# Move the selection/cursor the way the author's Excel session left it.
[void]$attrs.Range("F5").Select()

# New "Cytokine" sheet (lookup entity referenced from attribute row 6)
$cytokine = $wb.Worksheets.Add($null, $attrs)
$cytokine.Name = "Cytokine"
$cytokine.Range("A1").Value = "name"
$cytokine.Range("A2").Value = "IFNy"
$cytokine.Range("A3").Value = "IL17"
$cytokine.Range("A4").Value = "TNFA"
[void]$cytokine.Range("A4").Select()

# New "Stimulus" sheet (lookup entity referenced from attribute row 7)
$stimulus = $wb.Worksheets.Add($null, $cytokine)
$stimulus.Name = "Stimulus"
$stimulus.Range("A1").Value = "name"
$stimulus.Range("A2").Value = "E.Coli"
$stimulus.Range("A3").Value = "A.fumigatusconidia"
[void]$stimulus.Range("A4").Select()
